# "List of parts updated" — change every part whose STATUS is "Ordered"
# to "Ready" on the "Main" table (Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # column C = STATUS
    if ($cell.Text -eq "Ordered") {
        $cell.Value = "Ready"
    }
}
